$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.736.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.890.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.28'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4738'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2925'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06534'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.12'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07803'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.86'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.888.71'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7379'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.246'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.77'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.718.43'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.26'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007551'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.135.80'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.323'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.255'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.221'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.92'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.99'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.919'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.50%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09744'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.495'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.190'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04863'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.126'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6970'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01893'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.807'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.11'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.327'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.999'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4273'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8350'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.68'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.540'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.053'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.60'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '917.09'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05759'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.07%  '
